$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 44
$ws.Range("H44").Value = 49950
$ws.Range("J44").Value = 49950
$ws.Range("L44").Value = 49950
$ws.Range("N44").Value = -50874

# row 98
$ws.Range("H98").Value = 3155.1177
$ws.Range("I98").Value = 3170.6875
$ws.Range("J98").Value = 2906
$ws.Range("K98").Value = 3170.6875
$ws.Range("L98").Value = 2906
$ws.Range("M98").Value = -1672.6875
$ws.Range("N98").Value = -5902

# row 100
$ws.Range("H100").Value = 2500
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""

# row 122
$ws.Range("H122").Value = 3155.1177
$ws.Range("I122").Value = 3170.6875
$ws.Range("J122").Value = 2906
$ws.Range("K122").Value = 9512.0625
$ws.Range("L122").Value = 8718
$ws.Range("M122").Value = -7062.0625
$ws.Range("N122").Value = -13618

# row 135
$ws.Range("H135").Value = 1069.4546
$ws.Range("I135").Value = 972.9
$ws.Range("J135").Value = 2035
$ws.Range("K135").Value = 8756.1
$ws.Range("L135").Value = 18315
$ws.Range("M135").Value = -6221.1
$ws.Range("N135").Value = -23385

# row 138
$ws.Range("H138").Value = 3156.7036
$ws.Range("J138").Value = 3870.8809
$ws.Range("L138").Value = 11612.6427
$ws.Range("N138").Value = -21892.6427

# row 140
$ws.Range("H140").Value = 184997
$ws.Range("I140").Value = 169995
$ws.Range("J140").Value = 199999
$ws.Range("K140").Value = 169995
$ws.Range("L140").Value = 199999
$ws.Range("M140").Value = -164815
$ws.Range("N140").Value = -210359

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 829.1892
$ws.Range("J2").Value = 998.9167
$ws.Range("L2").Value = 998.9167
$ws.Range("N2").Value = -1224.9167

# row 24
$ws.Range("H24").Value = 100342.75
$ws.Range("J24").Value = 100342.75
$ws.Range("L24").Value = 100342.75
$ws.Range("N24").Value = -101090.75

# row 32
$ws.Range("H32").Value = 8339891.5
$ws.Range("I32").Value = 2602992.2
$ws.Range("K32").Value = 2602992.2
$ws.Range("M32").Value = -2602705.2

# row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

# row 74
$ws.Range("H74").Value = 2347.074
$ws.Range("I74").Value = 2064.6365
$ws.Range("K74").Value = 2064.6365
$ws.Range("M74").Value = -1190.6365

# row 77
$ws.Range("H77").Value = 2347.074
$ws.Range("I77").Value = 2064.6365
$ws.Range("K77").Value = 10323.1825
$ws.Range("M77").Value = -5955.182500000001

# row 98
$ws.Range("H98").Value = 33623.453
$ws.Range("J98").Value = 33623.453
$ws.Range("L98").Value = 33623.453
$ws.Range("N98").Value = -39613.453

# row 100
$ws.Range("H100").Value = 100342.75
$ws.Range("J100").Value = 100342.75
$ws.Range("L100").Value = 100342.75
$ws.Range("N100").Value = -102506.75

# row 116
$ws.Range("H116").Value = 829.1892
$ws.Range("J116").Value = 998.9167
$ws.Range("L116").Value = 998.9167
$ws.Range("N116").Value = -5586.9167

# row 122
$ws.Range("H122").Value = 4028.8367
$ws.Range("I122").Value = 2965.6453
$ws.Range("K122").Value = 8896.9359
$ws.Range("M122").Value = -6446.9359

# row 132
$ws.Range("H132").Value = 3651.9773
$ws.Range("I132").Value = 3133.5293
$ws.Range("J132").Value = 3978.4075
$ws.Range("K132").Value = 9400.5879
$ws.Range("L132").Value = 11935.2225
$ws.Range("M132").Value = -6870.5879
$ws.Range("N132").Value = -16995.2225

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 829.1892
$ws.Range("J3").Value = 998.9167
$ws.Range("L3").Value = 998.9167
$ws.Range("N3").Value = -1226.9167

# row 105
$ws.Range("H105").Value = 2597.1333
$ws.Range("I105").Value = 2915.08
$ws.Range("K105").Value = 2915.08
$ws.Range("M105").Value = -1168.08

$ws = $wb.Worksheets.Item("CRP")
# row 36
$ws.Range("H36").Value = 7829.8335
$ws.Range("I36").Value = 11666.667
$ws.Range("K36").Value = 11666.667
$ws.Range("M36").Value = -11278.667

# row 40
$ws.Range("H40").Value = 7829.8335
$ws.Range("I40").Value = 11666.667
$ws.Range("K40").Value = 11666.667
$ws.Range("M40").Value = -11506.667

# row 43
$ws.Range("H43").Value = 35415.668
$ws.Range("J43").Value = 35717.625
$ws.Range("L43").Value = 35717.625
$ws.Range("N43").Value = -36085.625

# row 75
$ws.Range("H75").Value = 69499
$ws.Range("J75").Value = 69499
$ws.Range("L75").Value = 69499
$ws.Range("N75").Value = -71495

# row 78
$ws.Range("H78").Value = 69499
$ws.Range("J78").Value = 69499
$ws.Range("L78").Value = 208497
$ws.Range("N78").Value = -218481

# row 95
$ws.Range("H95").Value = 22712.572
$ws.Range("J95").Value = 22712.572
$ws.Range("L95").Value = 22712.572
$ws.Range("N95").Value = -28204.572

# row 101
$ws.Range("H101").Value = 35415.668
$ws.Range("J101").Value = 35717.625
$ws.Range("L101").Value = 35717.625
$ws.Range("N101").Value = -42207.625

# row 122
$ws.Range("H122").Value = 1854.875
$ws.Range("I122").Value = 1854.875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5564.625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3114.625
$ws.Range("N122").Value = ""

# row 132
$ws.Range("H132").Value = 1198.25
$ws.Range("I132").Value = 1198.25
$ws.Range("K132").Value = 3594.75
$ws.Range("M132").Value = -1064.75

# row 141
$ws.Range("H141").Value = 2391999.2
$ws.Range("J141").Value = 2391999.2
$ws.Range("L141").Value = 2391999.2
$ws.Range("N141").Value = -2402359.2

$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 78146750
$ws.Range("I4").Value = 55062496
$ws.Range("J4").Value = 157892370
$ws.Range("K4").Value = 165187488
$ws.Range("L4").Value = 473677110
$ws.Range("M4").Value = -165187376
$ws.Range("N4").Value = -473677334

# row 39
$ws.Range("H39").Value = 4398.7334
$ws.Range("J39").Value = 4398.7334
$ws.Range("L39").Value = 13196.2002
$ws.Range("N39").Value = -13784.2002

# row 113
$ws.Range("H113").Value = 1221.2162
$ws.Range("J113").Value = 1191.2963
$ws.Range("L113").Value = 3573.8889
$ws.Range("N113").Value = -7913.8889

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 2614.3076
$ws.Range("I80").Value = 1749.75
$ws.Range("J80").Value = 2998.5557
$ws.Range("K80").Value = 1749.75
$ws.Range("L80").Value = 2998.5557
$ws.Range("M80").Value = -751.75
$ws.Range("N80").Value = -4994.5557

# row 83
$ws.Range("H83").Value = 2614.3076
$ws.Range("I83").Value = 1749.75
$ws.Range("J83").Value = 2998.5557
$ws.Range("K83").Value = 8748.75
$ws.Range("L83").Value = 14992.7785
$ws.Range("M83").Value = -3756.75
$ws.Range("N83").Value = -24976.7785

# row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = ""
$ws.Range("N92").Value = ""

# row 126
$ws.Range("H126").Value = 1554.0588
$ws.Range("J126").Value = 1396
$ws.Range("L126").Value = 4188
$ws.Range("N126").Value = -9128

# row 132
$ws.Range("H132").Value = 2113.8235
$ws.Range("I132").Value = 1075.2
$ws.Range("K132").Value = 3225.6
$ws.Range("M132").Value = -695.6000000000004

$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 1913.4546
$ws.Range("J46").Value = 5181.5
$ws.Range("L46").Value = 5181.5
$ws.Range("N46").Value = -5557.5

# row 82
$ws.Range("H82").Value = 4101.3125
$ws.Range("I82").Value = 2537.7144
$ws.Range("J82").Value = 5317.4443
$ws.Range("K82").Value = 2537.7144
$ws.Range("L82").Value = 5317.4443
$ws.Range("M82").Value = -2176.7144
$ws.Range("N82").Value = -6039.4443

# row 85
$ws.Range("H85").Value = 4101.3125
$ws.Range("I85").Value = 2537.7144
$ws.Range("J85").Value = 5317.4443
$ws.Range("K85").Value = 2537.7144
$ws.Range("L85").Value = 5317.4443
$ws.Range("M85").Value = -1289.7144
$ws.Range("N85").Value = -7813.4443

# row 95
$ws.Range("H95").Value = 39999
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").Value = ""

# row 103
$ws.Range("H103").Value = 4699.4
$ws.Range("J103").Value = 4699.4
$ws.Range("L103").Value = 4699.4
$ws.Range("N103").Value = -7043.4

# row 110
$ws.Range("H110").Value = 65847
$ws.Range("J110").Value = 65847
$ws.Range("L110").Value = 65847
$ws.Range("N110").Value = -74027

# row 122
$ws.Range("H122").Value = 4929424
$ws.Range("I122").Value = 17243882
$ws.Range("J122").Value = 3641
$ws.Range("K122").Value = 51731646
$ws.Range("L122").Value = 10923
$ws.Range("M122").Value = -51729196
$ws.Range("N122").Value = -15823

# row 132
$ws.Range("H132").Value = 3839.3125
$ws.Range("I132").Value = 2762.4546
$ws.Range("K132").Value = 8287.363799999999
$ws.Range("M132").Value = -5757.363799999999

$ws = $wb.Worksheets.Item("WVR")
# row 97
$ws.Range("H97").Value = 111000
$ws.Range("J97").Value = 111000
$ws.Range("L97").Value = 111000
$ws.Range("N97").Value = -112982

# row 107
$ws.Range("H107").Value = 571.7692
$ws.Range("I107").Value = 523.4
$ws.Range("J107").Value = 733
$ws.Range("K107").Value = 1570.2
$ws.Range("L107").Value = 2199
$ws.Range("M107").Value = 349.8000000000002
$ws.Range("N107").Value = -6039

# row 122
$ws.Range("H122").Value = 2016.7241
$ws.Range("I122").Value = 1929.1852
$ws.Range("J122").Value = 3198.5
$ws.Range("K122").Value = 5787.5556
$ws.Range("L122").Value = 9595.5
$ws.Range("M122").Value = -3337.5556
$ws.Range("N122").Value = -14495.5

# row 132
$ws.Range("H132").Value = 3247.6943
$ws.Range("I132").Value = 2790.7144
$ws.Range("K132").Value = 8372.143199999999
$ws.Range("M132").Value = -5842.143199999999

# row 137
$ws.Range("H137").Value = 74994
$ws.Range("J137").Value = 74994
$ws.Range("L137").Value = 74994
$ws.Range("N137").Value = -85194
